$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "weight" column header in J1, matching the style of the other headers
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "weight"

# Fill weight column (J2:J16) with value 1 for every data row
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 10).Value = 1
}

# Select whole column J (as would happen after inserting/selecting the new column)
$ws.Range("J1").Activate()
$ws.Columns.Item(10).Select()
